$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# CP003 -> CP003_Eminent (con "Sucursales Exclusivas Éminent:" como Dato1)
$ws.Range("B4").Value = "Sucursales Exclusivas Éminent:"
$ws.Range("A4").Value = "CP003_Eminent"

# Actualiza la seleccion activa de la hoja
$ws.Range("A4").Select()
